$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date serial number (45203 = 2023-10-04)
# for every data row (2 through 458). Bump it by one day to 45204
# (2023-10-05) to reflect the updated "last changed" timestamp.
$ws.Range("C2:C458").Value = 45204
